$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update results for Steel
$ws.Range("B3").Value = 1465.462811187185
$ws.Range("C4").Value = 33.93168714694505
$ws.Range("C5").Value = 2424.883513345701
$ws.Range("D8").Value = 126.8499407842673
